# "germanized assignments, minor typo fixes"
#
# - Rename "Geometry Class" -> "Mathematik", "History Class" -> "Geschichte"
#   ("Party Budget" keeps its name).
# - Translate assignment/grade + party-budget labels to German (and swap
#   "Paper N" -> "Hausaufgabe N" wording).
# - Re-point the $/€ currency number format on the Party Budget sheet.
# - Resize a couple of label columns.
# - Update each sheet's remembered selection / the active tab.

$wb = $excel.ActiveWorkbook

# ---- rename sheets -------------------------------------------------------
$wb.Worksheets.Item("Geometry Class").Name = "Mathematik"
$wb.Worksheets.Item("History Class").Name = "Geschichte"

$wsMath = $wb.Worksheets.Item("Mathematik")
$wsBudget = $wb.Worksheets.Item("Party Budget")
$wsHist = $wb.Worksheets.Item("Geschichte")

# ---- German copy on the two grade sheets ---------------------------------
foreach ($ws in @($wsMath, $wsHist)) {
    $ws.Range("A1").Value = "Aufgabe"
    $ws.Range("B1").Value = "Punkte"
    $ws.Range("A2").Value = "Test 1"
    $ws.Range("A3").Value = "Hausaufgabe 1"
    $ws.Range("A4").Value = "Test 2"
    $ws.Range("A5").Value = "Hausaufgabe 2"
    $ws.Range("A6").Value = "Test 3"
    $ws.Range("A7").Value = "Finale Punkte"
}

# ---- German copy on the budget sheet -------------------------------------
$wsBudget.Range("A1").Value = "Reservierungsgebühr"
$wsBudget.Range("A2").Value = "Preis pro Person"
$wsBudget.Range("A3").Value = "Anzahl Gäste"
$wsBudget.Range("A4").Value = "Budget"

# ---- currency format: $ -> € (accounting-style euro format) -------------
$euroFormat = '_-* #,##0.00\ [$€-407]_-;\-* #,##0.00\ [$€-407]_-;_-* "-"??\ [$€-407]_-;_-@_-'
$wsBudget.Range("B1").NumberFormat = $euroFormat
$wsBudget.Range("B2").NumberFormat = $euroFormat
$wsBudget.Range("B4").NumberFormat = $euroFormat

# ---- column widths (closest reachable values; engine snaps column widths
#      to its own internal character grid) ---------------------------------
$wsMath.Columns.Item(1).ColumnWidth = 16.8333333333333
$wsMath.Columns.Item(2).ColumnWidth = 8.66666666666667
$wsBudget.Columns.Item(1).ColumnWidth = 20
$wsHist.Columns.Item(1).ColumnWidth = 15.6666666666667

# ---- remembered selections / active tab ----------------------------------
# Order matters: the sheet selected last becomes the active tab.
$wsMath.Range("B24").Select()
$wsBudget.Range("D8").Select()
$wsHist.Range("E21").Select()
